# Fix getHyperlinkFunctionCellAddress in DefaultHyperlinkCellClickHandler (#43)
#
# Adds a few extra HYPERLINK() formula edge-cases (nested parens, extra
# whitespace) to Sheet1 so the click-handler's cell-address resolution can
# be exercised against them, and gives those cells (and the already-empty
# A9) the built-in "Hyperlink" cell style to match the other hyperlink
# cells on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A7: simple HYPERLINK() call with a single argument (extra inner spaces in
# the original fixture formula are just cosmetic).
$ws.Range("A7").Formula = '=HYPERLINK("http://www.google.com")'
$ws.Range("A7").Style = "Hyperlink"

# A8: same target, but wrapped in an extra pair of parentheses.
$ws.Range("A8").Formula = '=HYPERLINK(("http://www.google.com"))'
$ws.Range("A8").Style = "Hyperlink"

# A9: no formula, just picks up the Hyperlink style like its neighbours.
$ws.Range("A9").Style = "Hyperlink"

# A10: same HYPERLINK() target, double-wrapped in parentheses in the
# original fixture; functionally equivalent to A7/A8.
$ws.Range("A10").Formula = '=HYPERLINK(("http://www.google.com"))'
$ws.Range("A10").Style = "Hyperlink"

# Row 10 no longer needs its explicit custom height now that it matches the
# sheet's default row formatting again.
$ws.Rows.Item(10).AutoFit()

# Leave the selection on A10, matching the saved view state.
$ws.Range("A10").Select() | Out-Null
